# Inclusão de número de tcc e dados de orientador
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns ---
# New column C: "SIAPE" (between Orientador and CoOrientador)
$ws.Columns.Item(3).Insert()
# New column E: "NroTCC" (between CoOrientador and TítuloTCC)
$ws.Columns.Item(5).Insert()

# Approximate the original column widths for the two new columns (best effort;
# COM ColumnWidth only supports coarse granularity in this host).
$ws.Columns.Item(3).ColumnWidth = 9
$ws.Columns.Item(5).ColumnWidth = 13

# --- Header row (row 1) ---
$ws.Range("C1").Value = "SIAPE"
$ws.Range("E1").Value = "NroTCC"

# --- Data row (row 2) ---
# Columns that shifted but keep the placeholder-replacement text:
$ws.Range("A2").Value = "Aluno que está defendendo"
$ws.Range("B2").Value = "Orientador ABC"
$ws.Range("C2").Value = 12334443
$ws.Range("D2").Value = "Sou coorientador"
$ws.Range("E2").Value = 62
$ws.Range("F2").Value = "IFS TUTOR: UMA PROPOSTA DE CHATBOT PARA APOIO AOS ESTUDANTES DE NÍVEL SUPERIOR DO INSTITUTO FEDERAL DE SERGIPE"

$ws.Range("Q2").Value = "Professor da banca 1"
$ws.Range("S2").Value = "Professor da banca 2"
$ws.Range("U2").Value = "Prof coordenador"
$ws.Range("W2").Value = 9.7

# --- Selection / view state ---
$ws.Range("A3").Select()
